$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting old row 9 (totals) to 10 and old row 10 (footer) to 11
$ws.Rows("9:9").Insert()

# Copy formatting of row 8 into new row 9 (cell-by-cell to avoid merged-range copy quirks)
for ($col = 1; $col -le 17; $col++) {
  $src = $ws.Cells.Item(8, $col)
  $dst = $ws.Cells.Item(9, $col)
  $src.Copy()
  $dst.PasteSpecial(-4122)
}

# Re-apply merges for the new row 9 (mirrors rows 7/8 layout)
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# Fill in row 9 values (item #3: VOLTAREN)
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "VOLTAREN 75MG/3ML 3 AMP."
$ws.Range("H9").Value = "5:2"

# L9 and P9 must stay text even though their number format looks numeric;
# stage through a Text-formatted helper cell and paste values only.
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "1"
$ws.Range("ZZ1").Copy()
$ws.Range("L9").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("N9").Value = "51.00"

$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "16.8300"
$ws.Range("ZZ1").Copy()
$ws.Range("P9").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("Q9").Value = "0:1"

# Update the totals row (now row 10): old total 185 + new item's 16.83 = 201.83
$ws.Range("P10").Value = 201.83

# Update footer timestamp text (now row 11)
$ws.Range("A11").Value = "Monday, 1 September, 2025 9:48 AM"
